$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 216.96297
$ws.Range("I33").Value = 112.36364
$ws.Range("J33").Value = 677.2
$ws.Range("K33").Value = 112.36364
$ws.Range("L33").Value = 677.2
$ws.Range("M33").Value = 116.63636
$ws.Range("N33").Value = -1135.2
$ws.Range("H39").Value = 1216
$ws.Range("I39").Value = 619.7857
$ws.Range("K39").Value = 1859.3571
$ws.Range("M39").Value = -1563.3571
$ws.Range("H86").Value = 2134.0833
$ws.Range("I86").Value = 1198.75
$ws.Range("K86").Value = 1198.75
$ws.Range("M86").Value = -75.75
$ws.Range("H89").Value = 2134.0833
$ws.Range("I89").Value = 1198.75
$ws.Range("K89").Value = 5993.75
$ws.Range("M89").Value = -377.75
$ws.Range("H96").Value = 1161.1765
$ws.Range("I96").Value = 691.53845
$ws.Range("K96").Value = 2074.61535
$ws.Range("M96").Value = -701.61535
$ws.Range("H125").Value = 7145035
$ws.Range("I125").Value = 2640.75
$ws.Range("K125").Value = 23766.75
$ws.Range("M125").Value = -21306.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3705266.2
$ws.Range("I61").Value = 4763343
$ws.Range("J61").Value = 1997.5
$ws.Range("K61").Value = 4763343
$ws.Range("L61").Value = 1997.5
$ws.Range("M61").Value = -4763131
$ws.Range("N61").Value = -2421.5
$ws.Range("H110").Value = 2058.1667
$ws.Range("J110").Value = 2750
$ws.Range("L110").Value = 2750
$ws.Range("N110").Value = -6840
$ws.Range("H122").Value = 3551.611
$ws.Range("I122").Value = 3395.3333
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 10185.9999
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -7735.999899999999
$ws.Range("N122").Value = -17899
$ws.Range("H136").Value = 3705266.2
$ws.Range("I136").Value = 4763343
$ws.Range("J136").Value = 1997.5
$ws.Range("K136").Value = 14290029
$ws.Range("L136").Value = 5992.5
$ws.Range("M136").Value = -14287479
$ws.Range("N136").Value = -11092.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4455.4546
$ws.Range("I5").Value = 205
$ws.Range("J5").Value = 7997.5
$ws.Range("K5").Value = 205
$ws.Range("L5").Value = 7997.5
$ws.Range("M5").Value = -92
$ws.Range("N5").Value = -8223.5
$ws.Range("H26").Value = 33000
$ws.Range("I26").Value = 33000
$ws.Range("K26").Value = 33000
$ws.Range("M26").Value = -32708
$ws.Range("H96").Value = 21428
$ws.Range("I96").Value = 21428
$ws.Range("K96").Value = 21428
$ws.Range("M96").Value = -18682
$ws.Range("H105").Value = 2556.7144
$ws.Range("I105").Value = 2569.5
$ws.Range("J105").Value = 2524.75
$ws.Range("K105").Value = 2569.5
$ws.Range("L105").Value = 2524.75
$ws.Range("M105").Value = -822.5
$ws.Range("N105").Value = -6018.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 323.13333
$ws.Range("I7").Value = 351.25
$ws.Range("J7").Value = 291
$ws.Range("K7").Value = 351.25
$ws.Range("L7").Value = 291
$ws.Range("M7").Value = -238.25
$ws.Range("N7").Value = -517
$ws.Range("H16").Value = 1231.6666
$ws.Range("I16").Value = 963.3333
$ws.Range("K16").Value = 963.3333
$ws.Range("M16").Value = -676.3333
$ws.Range("H22").Value = 445.125
$ws.Range("I22").Value = 393
$ws.Range("J22").Value = 497.25
$ws.Range("K22").Value = 393
$ws.Range("L22").Value = 497.25
$ws.Range("M22").Value = -43
$ws.Range("N22").Value = -1197.25
$ws.Range("H36").Value = 3849
$ws.Range("I36").Value = 3849
$ws.Range("K36").Value = 3849
$ws.Range("M36").Value = -3461
$ws.Range("H40").Value = 3849
$ws.Range("I40").Value = 3849
$ws.Range("K40").Value = 3849
$ws.Range("M40").Value = -3689
$ws.Range("H113").Value = 1231.6666
$ws.Range("I113").Value = 963.3333
$ws.Range("K113").Value = 963.3333
$ws.Range("M113").Value = 1206.6667
$ws.Range("H132").Value = 230985.86
$ws.Range("J132").Value = 631152
$ws.Range("L132").Value = 1893456
$ws.Range("N132").Value = -1898516
$ws.Range("H134").Value = 2751.2778
$ws.Range("I134").Value = 2872
$ws.Range("J134").Value = 2147.6667
$ws.Range("K134").Value = 8616
$ws.Range("L134").Value = 6443.000100000001
$ws.Range("M134").Value = -6081
$ws.Range("N134").Value = -11513.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9730.223
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30224
$ws.Range("H23").Value = 175.54546
$ws.Range("I23").Value = 130.5
$ws.Range("J23").Value = 201.28572
$ws.Range("K23").Value = 391.5
$ws.Range("L23").Value = 603.85716
$ws.Range("M23").Value = -156.5
$ws.Range("N23").Value = -1073.85716
$ws.Range("H107").Value = 1272
$ws.Range("J107").Value = 1622.8334
$ws.Range("L107").Value = 4868.5002
$ws.Range("N107").Value = -8708.5002
$ws.Range("H113").Value = 1146.881
$ws.Range("J113").Value = 1182.975
$ws.Range("L113").Value = 3548.925
$ws.Range("N113").Value = -7888.924999999999
$ws.Range("H131").Value = 11599.363
$ws.Range("I131").Value = 1085.5714
$ws.Range("K131").Value = 3256.7142
$ws.Range("M131").Value = 1783.2858

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 50000
$ws.Range("J5").Value = 50000
$ws.Range("L5").Value = 50000
$ws.Range("N5").Value = -50224
$ws.Range("H122").Value = 4061.3777
$ws.Range("I122").Value = 2336.457
$ws.Range("K122").Value = 7009.370999999999
$ws.Range("M122").Value = -4559.370999999999
$ws.Range("H132").Value = 41671924
$ws.Range("I132").Value = 50003996
$ws.Range("K132").Value = 150011988
$ws.Range("M132").Value = -150009458

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 235.375
$ws.Range("I9").Value = 126.14286
$ws.Range("K9").Value = 126.14286
$ws.Range("M9").Value = 97.85714
$ws.Range("H22").Value = 1373.625
$ws.Range("I22").Value = 797.8
$ws.Range("J22").Value = 2333.3333
$ws.Range("K22").Value = 797.8
$ws.Range("L22").Value = 2333.3333
$ws.Range("M22").Value = -502.8
$ws.Range("N22").Value = -2923.3333
$ws.Range("H27").Value = 1373.625
$ws.Range("I27").Value = 797.8
$ws.Range("J27").Value = 2333.3333
$ws.Range("K27").Value = 797.8
$ws.Range("L27").Value = 2333.3333
$ws.Range("M27").Value = -690.8
$ws.Range("N27").Value = -2547.3333
$ws.Range("H30").Value = 505
$ws.Range("I30").Value = 505
$ws.Range("K30").Value = 505
$ws.Range("M30").Value = -397
$ws.Range("H35").Value = 1766.1428
$ws.Range("I35").Value = 1766.1428
$ws.Range("K35").Value = 1766.1428
$ws.Range("M35").Value = -1430.1428
$ws.Range("H122").Value = 4712.9536
$ws.Range("I122").Value = 4055.1562
$ws.Range("J122").Value = 6626.5454
$ws.Range("K122").Value = 12165.4686
$ws.Range("L122").Value = 19879.6362
$ws.Range("M122").Value = -9715.4686
$ws.Range("N122").Value = -24779.6362
$ws.Range("H136").Value = 5057.9165
$ws.Range("I136").Value = 2956.8572
$ws.Range("K136").Value = 8870.571599999999
$ws.Range("M136").Value = -6320.571599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5516.5
$ws.Range("I81").Value = 4366.6665
$ws.Range("J81").Value = 6666.3335
$ws.Range("K81").Value = 8733.333000000001
$ws.Range("L81").Value = 13332.667
$ws.Range("M81").Value = -7672.333000000001
$ws.Range("N81").Value = -15454.667
$ws.Range("H84").Value = 5516.5
$ws.Range("I84").Value = 4366.6665
$ws.Range("J84").Value = 6666.3335
$ws.Range("K84").Value = 43666.665
$ws.Range("L84").Value = 66663.33499999999
$ws.Range("M84").Value = -38362.665
$ws.Range("N84").Value = -77271.33499999999
$ws.Range("H100").Value = 1607.0834
$ws.Range("I100").Value = 1081
$ws.Range("J100").Value = 2659.25
$ws.Range("K100").Value = 2162
$ws.Range("L100").Value = 5318.5
$ws.Range("M100").Value = -1621
$ws.Range("N100").Value = -6400.5
$ws.Range("H107").Value = 2129.3572
$ws.Range("I107").Value = 1135.6666
$ws.Range("J107").Value = 3918
$ws.Range("K107").Value = 3406.9998
$ws.Range("L107").Value = 11754
$ws.Range("M107").Value = -1486.9998
$ws.Range("N107").Value = -15594
